$d = $word.ActiveDocument

$d.Content.Find.Execute("42×77=", $true, $false, $false, $false, $false, $true, 1, $false, "94×63=", 2)
$d.Content.Find.Execute("11×34=", $true, $false, $false, $false, $false, $true, 1, $false, "11×45=", 2)
$d.Content.Find.Execute("43×39=", $true, $false, $false, $false, $false, $true, 1, $false, "17×71=", 2)
$d.Content.Find.Execute("98×49=", $true, $false, $false, $false, $false, $true, 1, $false, "42×73=", 2)
$d.Content.Find.Execute("74×38=", $true, $false, $false, $false, $false, $true, 1, $false, "50×69=", 2)
$d.Content.Find.Execute("67×65=", $true, $false, $false, $false, $false, $true, 1, $false, "45×51=", 2)
$d.Content.Find.Execute("95×48=", $true, $false, $false, $false, $false, $true, 1, $false, "81×90=", 2)
$d.Content.Find.Execute("72×88=", $true, $false, $false, $false, $false, $true, 1, $false, "70×15=", 2)
$d.Content.Find.Execute("79×80=", $true, $false, $false, $false, $false, $true, 1, $false, "89×33=", 2)
$d.Content.Find.Execute("34×57=", $true, $false, $false, $false, $false, $true, 1, $false, "35×89=", 2)
$d.Content.Find.Execute("44×11=", $true, $false, $false, $false, $false, $true, 1, $false, "76×69=", 2)
$d.Content.Find.Execute("60×73=", $true, $false, $false, $false, $false, $true, 1, $false, "23×41=", 2)
$d.Content.Find.Execute("80×58=", $true, $false, $false, $false, $false, $true, 1, $false, "20×11=", 2)
$d.Content.Find.Execute("63×74=", $true, $false, $false, $false, $false, $true, 1, $false, "76×51=", 2)
$d.Content.Find.Execute("16×46=", $true, $false, $false, $false, $false, $true, 1, $false, "55×74=", 2)
$d.Content.Find.Execute("57×39=", $true, $false, $false, $false, $false, $true, 1, $false, "89×40=", 2)
$d.Content.Find.Execute("74×23=", $true, $false, $false, $false, $false, $true, 1, $false, "59×86=", 2)
$d.Content.Find.Execute("69×61=", $true, $false, $false, $false, $false, $true, 1, $false, "80×19=", 2)
$d.Content.Find.Execute("72×78=", $true, $false, $false, $false, $false, $true, 1, $false, "33×95=", 2)
$d.Content.Find.Execute("94×12=", $true, $false, $false, $false, $false, $true, 1, $false, "78×54=", 2)
$d.Content.Find.Execute("53×38=", $true, $false, $false, $false, $false, $true, 1, $false, "60×28=", 2)
$d.Content.Find.Execute("81×68=", $true, $false, $false, $false, $false, $true, 1, $false, "46×33=", 2)
$d.Content.Find.Execute("43×70=", $true, $false, $false, $false, $false, $true, 1, $false, "45×29=", 2)
$d.Content.Find.Execute("95×94=", $true, $false, $false, $false, $false, $true, 1, $false, "70×24=", 2)
$d.Content.Find.Execute("13×69=", $true, $false, $false, $false, $false, $true, 1, $false, "14×58=", 2)
